# Listeners class added to Runners.
# Appends 7 new test-result rows (rows 243-249) to the bottom of the
# existing data table on sheet1, reusing the existing shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("Login with valid username and password", "PASSED", "chrome", "09.11.22"),
    @("Login Failure with invalid username or password", "FAILED", "chrome", "09.11.22"),
    @("Login Failure with invalid username or password", "FAILED", "chrome", "09.11.22"),
    @("Login Failure with invalid username or password", "FAILED", "chrome", "09.11.22"),
    @("Login Failure with invalid username or password", "PASSED", "chrome", "09.11.22"),
    @("Login Failure with invalid username or password", "PASSED", "chrome", "09.11.22"),
    @("Login Failure with invalid username or password", "PASSED", "chrome", "09.11.22")
)

$startRow = 243
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
}
